$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells whose shared-string lookups changed ---
# Row 2: TESTNAME row - A/B swapped to new key/de text (C stays "Speech Rhythm..." in en)
$ws.Range("A2").Value = "TESTNAME"
$ws.Range("B2").Value = "Test zur Rhythmuswahrnehmung in der Sprache"

# Row 3: INSTRUCTIONS row - de/en long text bodies
$ws.Range("B3").Value = "Wir werden Ihnen nun einige Reihen von jeweils vier einzelnen Wörtern oder kurzen Ausdrücken präsentieren, z.B. Teller – Schachtel – Ratte – Parkett. Bitte lesen Sie diese Reihen aufmerksam und achten Sie dabei besonders auf den Sprachrhythmus, d.h. auf die Betonungsmuster der Wörter. <br/> In manchen dieser Reihen weisen alle Elemente, d.h. alle enthaltenen Einzelwörter oder Mehrwort-Ausdrücke, dasselbe Betonungsmuster auf.  In einigen dieser Reihen ist jedoch genau ein Element enthalten, das ein anderes Betonungsmuster aufweist als die übrigen – im Beispiel oben ist es das  Wort „Parkett“. Ihre Aufgabe besteht darin, für jede Reihe zu beurteilen, ob ein Element mit abweichendem Betonungsmuster enthalten ist oder nicht. Falls alle Elemente der Reihe das gleiche Betonungsmuster aufweisen, klicken Sie bitte auf das Feld „Alle gleich“. Falls ein Element ein abweichendes Betonungsmuster aufweist, klicken Sie bitte auf das Feld dieses abweichenden Elements.<br/> Bitte beachten Sie, dass Sie für jede Reihe **{{time_out}} Sekunden Zeit** haben, um Ihr Urteil abzugeben."
$ws.Range("C3").Value = 'We will now present you with some series of four individual words or short expressions each, e.g. plate - box - rat - parquet.  Please read these series carefully, paying particular  attention to the rhythm of speech, i.e. the stress patterns of the words. <br/> In some of these series, however, there is exactly one element that has a different stress pattern than the others - in the example above it is the word "parquet".  Your task is to judge for each row whether or not there is an element with a different stress pattern.  If all the elements in the row have the same stress pattern, please click on the field "All equal".  If an element has a different stress pattern, please click on the field of this different element. <br/>Please note that you have **{{time_out}} seconds** for each row to make your judgement.'

# Row 4: PROMPT row - de/en text (kept full / long version here; shortened version moves to new PROMPT_SHORT row)
$ws.Range("B4").Value = "Ist der Rhythmus aller Wörter gleich oder hat eines ein abweichendes Betonungsmuster?<br/>  Wählen Sie “Alle gleich” oder klicken Sie auf das Wort, das abweicht. <br/> Sie haben {{time_out}} Sekunden Zeit zu antworten.`n"
$ws.Range("C4").Value = 'Is the rhythm of all words equal or does one not fit to the others? Click on the word, that does not fit or on "All equal". <br/>You got {{time_out}} seconds to answer.'

# Row 6: FEEDBACK row
$ws.Range("A6").Value = "FEEDBACK"
$ws.Range("B6").Value = "Sie haben {{num_correct}} von {{num_items}} Fragen richtig beantwortet ({{perc_correct}}%)."
$ws.Range("C6").Value = "You answered {{num_correct}} out of {{num_items}} questions correct ({{perc_correct}}%)."

# Row 9: WELCOME row
$ws.Range("B9").Value = "Test: Rhythmuswahrnehmung in der Sprache"

# --- Rows 12-16: replace the old EXAMPLE/EXAMPLE_PROMPT/A Sample Question block
#     with the new EXAMPLE1/EXAMPLE_PROMPT1/EXAMPLE_FEEDBACK_CORRECT1 block ---
$ws.Range("A12").Value = "EXAMPLE1"
$ws.Range("B12").Value = "Erste Beispielaufgabe"
$ws.Range("C12").Value = "First Sample Question"

$ws.Range("A13").Value = "EXAMPLE_PROMPT1"
$ws.Range("B13").Value = "Hier zunächst ein Beispiel."
$ws.Range("C13").Value = "First, an example."

$ws.Range("A14").Value = "EXAMPLE_FEEDBACK_CORRECT1"
$ws.Range("B14").Value = "Korrekt, ReVIER hat ein abweichendes Betonungsmuster."
$ws.Range("C14").Value = "Correct, ReVIER did not fit rhythmically."

$ws.Range("A15").Value = "EXAMPLE_FEEDBACK_INCORRECT"
$ws.Range("B15").Value = "Das war leider nicht richtig. Versuchen Sie es nochmal!"
$ws.Range("C15").Value = "The answer was not correct. Please try again!"

$ws.Range("A16").Value = "EXAMPLE_FEEDBACK_TOO_SLOW"
$ws.Range("B16").Value = "Das war leider zu langsam. Versuchen Sie es nochmal!"
$ws.Range("C16").Value = "This was too slow. Please try again!"

# --- New rows 17-20: second sample item + shortened prompt ---
$ws.Range("A17").Value = "EXAMPLE_FEEDBACK_CORRECT2"
$ws.Range("B17").Value = "Korrekt, alle hatten dasselbe Betonungsmuster."
$ws.Range("C17").Value = "Correct, all had the same rhythm."

$ws.Range("A18").Value = "EXAMPLE2"
$ws.Range("B18").Value = "Zweite Beispielaufgabe"
$ws.Range("C18").Value = "Second Sample Question"

$ws.Range("A19").Value = "EXAMPLE_PROMPT2"
$ws.Range("B19").Value = "Hier noch ein Beispiel."
$ws.Range("C19").Value = "Another example."

$ws.Range("A20").Value = "PROMPT_SHORT"
$ws.Range("B20").Value = "Ist der Rhythmus aller Wörter gleich oder hat eines ein abweichendes Betonungsmuster?`n"
$ws.Range("C20").Value = "Is the rhythm of all words equal or does one not fit to the others?"

# --- Formatting: new rows need the same wrap/top-align style (s="1") as the rest of the table ---
$ws.Range("A17:C20").WrapText = $true
$ws.Range("A17:C20").VerticalAlignment = -4160

# --- Row heights (match taller wrapped rows) ---
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 60

# --- View / selection state ---
$ws.Range("A20:C20").Select()
